$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new DAMSLTag (column I / 9), new DialogAct (column J / 10)
$updates = @(
    @(29, 'ba', 'Appreciation'),
    @(40, 'ba', 'Appreciation'),
    @(41, 'sd', 'Statement-non-opinion'),
    @(51, 'sv', 'Statement-opinion'),
    @(55, 'sd', 'Statement-non-opinion'),
    @(65, 'sv', 'Statement-opinion'),
    @(70, 'sd', 'Statement-non-opinion'),
    @(73, 'sd', 'Statement-non-opinion'),
    @(80, 'sv', 'Statement-opinion'),
    @(86, 'sd', 'Statement-non-opinion'),
    @(90, 'qy', 'Yes-No-Question'),
    @(108, 'sd', 'Statement-non-opinion'),
    @(117, 'sd', 'Statement-non-opinion'),
    @(131, 'sd', 'Statement-non-opinion'),
    @(145, 'b', 'Acknowledge (Backchannel)'),
    @(149, 'sd', 'Statement-non-opinion'),
    @(161, 'sd', 'Statement-non-opinion'),
    @(167, 'sv', 'Statement-opinion'),
    @(170, 'sv', 'Statement-opinion'),
    @(172, 'ba', 'Appreciation'),
    @(181, 'ba', 'Appreciation'),
    @(190, 'sv', 'Statement-opinion'),
    @(192, 'ba', 'Appreciation'),
    @(193, 'sd', 'Statement-non-opinion'),
    @(194, 'aa', 'Agree/Accept'),
    @(211, 'sv', 'Statement-opinion'),
    @(232, 'b', 'Acknowledge (Backchannel)'),
    @(241, 'sv', 'Statement-opinion'),
    @(242, 'ba', 'Appreciation'),
    @(249, '%', 'Uninterpretable'),
    @(254, 'ba', 'Appreciation'),
    @(255, 'b', 'Acknowledge (Backchannel)'),
    @(257, 'sv', 'Statement-opinion'),
    @(261, 'aa', 'Agree/Accept'),
    @(280, 'sv', 'Statement-opinion'),
    @(282, 'sd', 'Statement-non-opinion'),
    @(283, 'sd', 'Statement-non-opinion'),
    @(286, 'sd', 'Statement-non-opinion'),
    @(295, 'aa', 'Agree/Accept'),
    @(306, 'sv', 'Statement-opinion'),
    @(315, 'b', 'Acknowledge (Backchannel)'),
    @(330, 'sv', 'Statement-opinion'),
    @(334, 'sd', 'Statement-non-opinion'),
    @(335, 'sd', 'Statement-non-opinion'),
    @(342, 'sd', 'Statement-non-opinion'),
    @(358, 'sv', 'Statement-opinion'),
    @(359, 'sv', 'Statement-opinion'),
    @(367, 'sd', 'Statement-non-opinion'),
    @(384, 'sd', 'Statement-non-opinion'),
    @(391, 'b', 'Acknowledge (Backchannel)'),
    @(407, 'aa', 'Agree/Accept'),
    @(408, 'ba', 'Appreciation'),
    @(423, 'sv', 'Statement-opinion'),
    @(427, 'sv', 'Statement-opinion'),
    @(437, 'b', 'Acknowledge (Backchannel)'),
    @(438, 'ba', 'Appreciation'),
    @(443, 'sd', 'Statement-non-opinion'),
    @(447, 'ba', 'Appreciation'),
    @(460, 'sv', 'Statement-opinion'),
    @(474, 'sv', 'Statement-opinion'),
    @(479, 'sd', 'Statement-non-opinion'),
    @(491, 'sv', 'Statement-opinion'),
    @(496, '%', 'Uninterpretable'),
    @(508, 'sd', 'Statement-non-opinion'),
    @(510, 'sv', 'Statement-opinion'),
    @(516, 'b', 'Acknowledge (Backchannel)'),
    @(527, 'ba', 'Appreciation'),
    @(534, 'ba', 'Appreciation'),
    @(539, 'sd', 'Statement-non-opinion'),
    @(545, 'sd', 'Statement-non-opinion'),
    @(556, 'sv', 'Statement-opinion'),
    @(560, 'sv', 'Statement-opinion'),
    @(569, 'ba', 'Appreciation'),
    @(574, 'sd', 'Statement-non-opinion'),
    @(589, 'sd', 'Statement-non-opinion'),
    @(591, 'ba', 'Appreciation'),
    @(596, 'sv', 'Statement-opinion'),
    @(603, 'sv', 'Statement-opinion'),
)

foreach ($u in $updates) {
    $row = $u[0]
    $damslTag = $u[1]
    $dialogAct = $u[2]
    $ws.Cells.Item($row, 9).Value = $damslTag
    $ws.Cells.Item($row, 10).Value = $dialogAct
}
